$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the diagonal of newly-computed error values, continuing the
# staircase pattern already present in the sheet (one additional
# column per row, moving one column to the right as the row index
# increases).
$ws.Range("J16").Value = 0.1485032540344368
$ws.Range("I17").Value = 0.1734537503564907
$ws.Range("H18").Value = 0.2422520263583712
$ws.Range("G19").Value = 0.2534537503564908
$ws.Range("F20").Value = 0.381103329907261
$ws.Range("E21").Value = 0.04235042473292953
$ws.Range("D22").Value = 0.07961008106920435
$ws.Range("C23").Value = 0.02893023050567838
$ws.Range("B24").Value = 0.02940328597706714
